$p = $ppt.ActivePresentation

$titleText = "Onix Renewable Signs MoU with Government of Gujarat at Vibrant Gujarat 2026 for Clean Energy Expansion - Energetica India Magazine"

# --- Slide 1 ---
$s1 = $p.Slides.Item(1)
$shape1 = $s1.Shapes.Item(1)
$tf1 = $shape1.TextFrame

# Title
$tf1.TextRange.Paragraphs(1).Text = ""
$tf1.TextRange.Paragraphs(1).Text = $titleText
$tf1.TextRange.Paragraphs(1).Font.Color.RGB = 0xFFFFFF

# Bullet 1
$tf1.TextRange.Paragraphs(2).Text = ""
$tf1.TextRange.Paragraphs(2).Text = "Onix Renewable has signed a Memorandum of Understanding (MoU) with the Government of Gujarat."
$tf1.TextRange.Paragraphs(2).Font.Color.RGB = 0xF0F0F0

# Bullet 2
$tf1.TextRange.Paragraphs(3).Text = ""
$tf1.TextRange.Paragraphs(3).Text = "The signing took place during the Vibrant Gujarat 2026 event."
$tf1.TextRange.Paragraphs(3).Font.Color.RGB = 0xF0F0F0

# Position / size (points chosen so float32 truncation yields exact target EMU)
$shape1.Left = 57.600001
$shape1.Top = 57.600001
$shape1.Width = 604.8000189999999
$shape1.Height = 179.999993

# --- Slide 2 ---
$s2 = $p.Slides.Item(2)
$shape2 = $s2.Shapes.Item(1)
$tf2 = $shape2.TextFrame

# Title
$tf2.TextRange.Paragraphs(1).Text = ""
$tf2.TextRange.Paragraphs(1).Text = $titleText
$tf2.TextRange.Paragraphs(1).Font.Color.RGB = 0xFFFFFF

# Bullet 1
$tf2.TextRange.Paragraphs(2).Text = ""
$tf2.TextRange.Paragraphs(2).Text = "The MoU aims to facilitate clean energy expansion in the region."
$tf2.TextRange.Paragraphs(2).Font.Color.RGB = 0xF0F0F0

# Bullet 2
$tf2.TextRange.Paragraphs(3).Text = ""
$tf2.TextRange.Paragraphs(3).Text = "This partnership is part of Gujarat's broader initiative to enhance renewable energy capacity."
$tf2.TextRange.Paragraphs(3).Font.Color.RGB = 0xF0F0F0

# Position / size (points chosen so float32 truncation yields exact target EMU)
$shape2.Left = 57.600001
$shape2.Top = 57.600001
$shape2.Width = 604.8000189999999
$shape2.Height = 179.999993

# --- Slide 3 ---
$s3 = $p.Slides.Item(3)
$shape3 = $s3.Shapes.Item(1)
$tf3 = $shape3.TextFrame

# Title
$tf3.TextRange.Paragraphs(1).Text = ""
$tf3.TextRange.Paragraphs(1).Text = $titleText
$tf3.TextRange.Paragraphs(1).Font.Color.RGB = 0xFFFFFF

# Bullet 1
$tf3.TextRange.Paragraphs(2).Text = ""
$tf3.TextRange.Paragraphs(2).Text = "The agreement aligns with India's commitment to increasing renewable energy sources."
$tf3.TextRange.Paragraphs(2).Font.Color.RGB = 0xF0F0F0

# Bullet 2
$tf3.TextRange.Paragraphs(3).Text = ""
$tf3.TextRange.Paragraphs(3).Text = "Gujarat is known for its significant investments in clean energy projects."
$tf3.TextRange.Paragraphs(3).Font.Color.RGB = 0xF0F0F0

# Position / size (points chosen so float32 truncation yields exact target EMU)
$shape3.Left = 57.600001
$shape3.Top = 57.600001
$shape3.Width = 604.8000189999999
$shape3.Height = 179.999993

# --- Slide 4 ---
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(1)
$tf4 = $shape4.TextFrame

# Title
$tf4.TextRange.Paragraphs(1).Text = ""
$tf4.TextRange.Paragraphs(1).Text = $titleText
$tf4.TextRange.Paragraphs(1).Font.Color.RGB = 0xFFFFFF

# Bullet 1
$tf4.TextRange.Paragraphs(2).Text = ""
$tf4.TextRange.Paragraphs(2).Text = "The MoU is expected to contribute to job creation in the clean energy sector."
$tf4.TextRange.Paragraphs(2).Font.Color.RGB = 0xF0F0F0

# Bullet 2
$tf4.TextRange.Paragraphs(3).Text = ""
$tf4.TextRange.Paragraphs(3).Text = "Onix Renewable aims to leverage Gujarat's favorable policies for renewable energy."
$tf4.TextRange.Paragraphs(3).Font.Color.RGB = 0xF0F0F0

# Position / size (points chosen so float32 truncation yields exact target EMU)
$shape4.Left = 57.600001
$shape4.Top = 57.600001
$shape4.Width = 604.8000189999999
$shape4.Height = 179.999993

# --- Slide 5 ---
$s5 = $p.Slides.Item(5)
$shape5 = $s5.Shapes.Item(1)
$tf5 = $shape5.TextFrame

# Title
$tf5.TextRange.Paragraphs(1).Text = ""
$tf5.TextRange.Paragraphs(1).Text = $titleText
$tf5.TextRange.Paragraphs(1).Font.Color.RGB = 0xFFFFFF

# Bullet 1
$tf5.TextRange.Paragraphs(2).Text = ""
$tf5.TextRange.Paragraphs(2).Text = "The Vibrant Gujarat summit is a platform for investment and collaboration in various sectors."
$tf5.TextRange.Paragraphs(2).Font.Color.RGB = 0xF0F0F0

# Bullet 2
$tf5.TextRange.Paragraphs(3).Text = ""
$tf5.TextRange.Paragraphs(3).Text = "Onix Renewable's initiative reflects the growing trend of public-private partnerships in clean energy."
$tf5.TextRange.Paragraphs(3).Font.Color.RGB = 0xF0F0F0

# Position / size (points chosen so float32 truncation yields exact target EMU)
$shape5.Left = 57.600001
$shape5.Top = 57.600001
$shape5.Width = 604.8000189999999
$shape5.Height = 179.999993

# --- Slide 6 ---
$s6 = $p.Slides.Item(6)
$shape6 = $s6.Shapes.Item(1)
$tf6 = $shape6.TextFrame

# Title
$tf6.TextRange.Paragraphs(1).Text = ""
$tf6.TextRange.Paragraphs(1).Text = $titleText
$tf6.TextRange.Paragraphs(1).Font.Color.RGB = 0xFFFFFF

# Bullet 1
$tf6.TextRange.Paragraphs(2).Text = ""
$tf6.TextRange.Paragraphs(2).Text = "The MoU signifies a strategic move towards sustainable energy solutions."
$tf6.TextRange.Paragraphs(2).Font.Color.RGB = 0xF0F0F0

# Bullet 2
$tf6.TextRange.Paragraphs(3).Text = ""
$tf6.TextRange.Paragraphs(3).Text = "This collaboration is part of Gujarat's vision for a greener future."
$tf6.TextRange.Paragraphs(3).Font.Color.RGB = 0xF0F0F0

# Position / size (points chosen so float32 truncation yields exact target EMU)
$shape6.Left = 57.600001
$shape6.Top = 57.600001
$shape6.Width = 604.8000189999999
$shape6.Height = 179.999993
